# ECB statements for classification
# Fill in the "Classification Rui" column (G) for the ecb_monetary_policy_decisions sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ecb_monetary_policy_decisions")

$values = @(
    "Neutral",   # G2  - 2024-03-07
    "Neutral",   # G3  - 2024-01-25
    "Neutral",   # G4  - 2023-12-14
    "Neutral",   # G5  - 2023-10-26
    "Hawkish",   # G6  - 2023-09-14
    "Hawkish",   # G7  - 2023-07-27
    "Hawkish",   # G8  - 2023-06-15
    "Hawkish",   # G9  - 2023-05-04
    "Hawkish",   # G10 - 2023-03-16
    "Hawkish",   # G11 - 2023-02-02
    "Hawkish",   # G12 - 2022-12-15
    "Hawkish",   # G13 - 2022-10-27
    "Hawkish",   # G14 - 2022-09-08
    "Hawkish"    # G15 - 2022-07-21
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}

$wb.Save()
